# Apply "Update July Extended Experiment" changes to results sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "destination"
$ws.Range("D2").Value = "sum"
$ws.Range("E2").Value = 1.414213562373095

$ws.Range("B3").Value = "destination"
$ws.Range("C3").Value = "departuredelay"
$ws.Range("D3").Value = "sum"
$ws.Range("E3").Value = 1.414213562373095

$ws.Range("B4").Value = "destination"
$ws.Range("C4").Value = "weatherdelay"
$ws.Range("D4").Value = "sum"
$ws.Range("E4").Value = 1.414213562373095

$ws.Range("B5").Value = "destination"
$ws.Range("C5").Value = "distance"
$ws.Range("D5").Value = "sum"
$ws.Range("E5").Value = 1.414213562373095

$ws.Range("B6").Value = "destination"
$ws.Range("C6").Value = "arrivaldelay"
$ws.Range("D6").Value = "avg"
$ws.Range("E6").Value = 1.414213562373095

$ws.Range("B7").Value = "destination"
$ws.Range("C7").Value = "departuredelay"
$ws.Range("E7").Value = 1.414213562373095

$ws.Range("B8").Value = "destination"
$ws.Range("C8").Value = "weatherdelay"
$ws.Range("D8").Value = "avg"
$ws.Range("E8").Value = 1.414213562373095

$ws.Range("B9").Value = "destination"
$ws.Range("C9").Value = "distance"
$ws.Range("D9").Value = "avg"
$ws.Range("E9").Value = 1.414213562373095

$ws.Range("B10").Value = "destination"
$ws.Range("D10").Value = "max"
$ws.Range("E10").Value = 1.414213562373095

$ws.Range("B11").Value = "destination"
$ws.Range("D11").Value = "max"
$ws.Range("E11").Value = 1.414213562373095

$ws.Range("B12").Value = "destination"
$ws.Range("D12").Value = "max"
$ws.Range("E12").Value = 1.414213562373095

$ws.Range("B13").Value = "destination"
$ws.Range("D13").Value = "max"
$ws.Range("E13").Value = 1.414213562373095

$ws.Range("B14").Value = "carrier"
$ws.Range("C14").Value = "distance"
$ws.Range("D14").Value = "sum"
$ws.Range("E14").Value = 1.304700280674785

$ws.Range("B15").Value = "carrier"
$ws.Range("C15").Value = "arrivaldelay"
$ws.Range("D15").Value = "sum"
$ws.Range("E15").Value = 1.303992117540466

$ws.Range("B16").Value = "carrier"
$ws.Range("D16").Value = "sum"
$ws.Range("E16").Value = 1.287187434832667

$ws.Range("B17").Value = "carrier"
$ws.Range("C17").Value = "departuredelay"
$ws.Range("D17").Value = "sum"
$ws.Range("E17").Value = 1.284000589640841

$ws.Range("B18").Value = "origin"
$ws.Range("D18").Value = "sum"
$ws.Range("E18").Value = 1.181478964005186

$ws.Range("B19").Value = "origin"
$ws.Range("C19").Value = "departuredelay"
$ws.Range("E19").Value = 1.159562473740353

$ws.Range("C20").Value = "departuredelay"
$ws.Range("D20").Value = "sum"
$ws.Range("E20").Value = 1.149615081834924

$ws.Range("C21").Value = "departuredelay"
$ws.Range("E21").Value = 1.146319502273639

$ws.Range("B22").Value = "origin"
$ws.Range("C22").Value = "weatherdelay"
$ws.Range("D22").Value = "max"
$ws.Range("E22").Value = 1.137307698092018

$ws.Range("C23").Value = "weatherdelay"
$ws.Range("D23").Value = "sum"
$ws.Range("E23").Value = 1.132535687789139

$ws.Range("C24").Value = "distance"
$ws.Range("D24").Value = "avg"
$ws.Range("E24").Value = 1.123897733019126

$ws.Range("C25").Value = "weatherdelay"
$ws.Range("D25").Value = "avg"
$ws.Range("E25").Value = 1.110854669031078

$ws.Range("C26").Value = "arrivaldelay"
$ws.Range("E26").Value = 1.08948508648008

$ws.Range("C27").Value = "arrivaldelay"
$ws.Range("D27").Value = "avg"
$ws.Range("E27").Value = 1.087737404594133

$ws.Range("C28").Value = "distance"
$ws.Range("D28").Value = "max"
$ws.Range("E28").Value = 1.047477046884439

$ws.Range("C29").Value = "distance"
$ws.Range("D29").Value = "sum"
$ws.Range("E29").Value = 0.9968480213809268

$ws.Range("B30").Value = "carrier"
$ws.Range("E30").Value = 0.9399178319057042

$ws.Range("B31").Value = "carrier"
$ws.Range("E31").Value = 0.8599934570464456

$ws.Range("B32").Value = "carrier"
$ws.Range("C32").Value = "weatherdelay"
$ws.Range("D32").Value = "avg"
$ws.Range("E32").Value = 0.7908914482216012

$ws.Range("B33").Value = "carrier"
$ws.Range("C33").Value = "departuredelay"
$ws.Range("D33").Value = "avg"
$ws.Range("E33").Value = 0.7879647717808693

$ws.Range("B34").Value = "carrier"
$ws.Range("D34").Value = "avg"
$ws.Range("E34").Value = 0.7377034608157825

$ws.Range("B35").Value = "year"
$ws.Range("C35").Value = "departuredelay"
$ws.Range("D35").Value = "max"
$ws.Range("E35").Value = 0.7362206236615845

$ws.Range("B36").Value = "carrier"
$ws.Range("C36").Value = "distance"
$ws.Range("D36").Value = "max"
$ws.Range("E36").Value = 0.7268269685757737

$ws.Range("B37").Value = "carrier"
$ws.Range("D37").Value = "avg"
$ws.Range("E37").Value = 0.7127463816166288

$ws.Range("B38").Value = "year"
$ws.Range("C38").Value = "arrivaldelay"
$ws.Range("D38").Value = "sum"
$ws.Range("E38").Value = 0.7038647382433287

$ws.Range("C39").Value = "departuredelay"
$ws.Range("E39").Value = 0.6986967117610029

$ws.Range("C40").Value = "weatherdelay"
$ws.Range("D40").Value = "sum"
$ws.Range("E40").Value = 0.6938926425140296

$ws.Range("B41").Value = "year"
$ws.Range("D41").Value = "sum"
$ws.Range("E41").Value = 0.6786621690247974

$ws.Range("B42").Value = "carrier"
$ws.Range("C42").Value = "arrivaldelay"
$ws.Range("E42").Value = 0.6767445554101715

$ws.Range("B43").Value = "year"
$ws.Range("D43").Value = "sum"
$ws.Range("E43").Value = 0.6610645547736658

$ws.Range("B44").Value = "year"
$ws.Range("E44").Value = 0.5917881802992037

$ws.Range("B45").Value = "month"
$ws.Range("C45").Value = "weatherdelay"
$ws.Range("D45").Value = "max"
$ws.Range("E45").Value = 0.5686392777393928

$ws.Range("B46").Value = "week"
$ws.Range("E46").Value = 0.5675596982230317

$ws.Range("B47").Value = "week"
$ws.Range("C47").Value = "departuredelay"
$ws.Range("E47").Value = 0.5578387397648435

$ws.Range("B48").Value = "day"
$ws.Range("C48").Value = "weatherdelay"
$ws.Range("E48").Value = 0.5538619379975623

$ws.Range("B49").Value = "month"
$ws.Range("C49").Value = "departuredelay"
$ws.Range("D49").Value = "max"
$ws.Range("E49").Value = 0.5479432657805892

$ws.Range("B50").Value = "month"
$ws.Range("C50").Value = "weatherdelay"
$ws.Range("E50").Value = 0.4255643779824537

$ws.Range("B51").Value = "month"
$ws.Range("E51").Value = 0.3985042978667666

$ws.Range("B52").Value = "year"
$ws.Range("C52").Value = "departuredelay"
$ws.Range("D52").Value = "avg"
$ws.Range("E52").Value = 0.3953204047461664

$ws.Range("B53").Value = "day"
$ws.Range("D53").Value = "sum"
$ws.Range("E53").Value = 0.3751430180375182

$ws.Range("B54").Value = "month"
$ws.Range("C54").Value = "weatherdelay"
$ws.Range("E54").Value = 0.3620494216021585

$ws.Range("E55").Value = 0.3471053241921532

$ws.Range("B56").Value = "day"
$ws.Range("C56").Value = "departuredelay"
$ws.Range("E56").Value = 0.345219412202222

$ws.Range("B57").Value = "day"
$ws.Range("C57").Value = "weatherdelay"
$ws.Range("D57").Value = "sum"
$ws.Range("E57").Value = 0.3446703710278781

$ws.Range("B58").Value = "month"
$ws.Range("D58").Value = "avg"
$ws.Range("E58").Value = 0.3397189862389229

$ws.Range("B59").Value = "year"
$ws.Range("C59").Value = "distance"
$ws.Range("D59").Value = "max"
$ws.Range("E59").Value = 0.3362550616772066

$ws.Range("B60").Value = "day"
$ws.Range("C60").Value = "weatherdelay"
$ws.Range("D60").Value = "avg"
$ws.Range("E60").Value = 0.3105376281614586

$ws.Range("B61").Value = "week"
$ws.Range("C61").Value = "distance"
$ws.Range("D61").Value = "max"
$ws.Range("E61").Value = 0.3082372977329072

$ws.Range("C62").Value = "departuredelay"
$ws.Range("E62").Value = 0.3079001004573763

$ws.Range("B63").Value = "year"
$ws.Range("C63").Value = "distance"
$ws.Range("E63").Value = 0.3078557825256665

$ws.Range("B64").Value = "week"
$ws.Range("C64").Value = "weatherdelay"
$ws.Range("E64").Value = 0.2960129252291743

$ws.Range("C65").Value = "distance"
$ws.Range("E65").Value = 0.2828235476907014

$ws.Range("C66").Value = "arrivaldelay"
$ws.Range("E66").Value = 0.2698660745389153

$ws.Range("C67").Value = "distance"
$ws.Range("D67").Value = "max"
$ws.Range("E67").Value = 0.2681403376901599

$ws.Range("B68").Value = "day"
$ws.Range("C68").Value = "distance"
$ws.Range("D68").Value = "sum"
$ws.Range("E68").Value = 0.2673315557244986

$ws.Range("B69").Value = "day"
$ws.Range("C69").Value = "arrivaldelay"
$ws.Range("D69").Value = "sum"
$ws.Range("E69").Value = 0.2591010417324213

$ws.Range("B70").Value = "week"
$ws.Range("E70").Value = 0.2528766172066888

$ws.Range("B71").Value = "day"
$ws.Range("C71").Value = "distance"
$ws.Range("D71").Value = "max"
$ws.Range("E71").Value = 0.2438694113483663

$ws.Range("B72").Value = "week"
$ws.Range("C72").Value = "weatherdelay"
$ws.Range("D72").Value = "avg"
$ws.Range("E72").Value = 0.2334715372373947

$ws.Range("E73").Value = 0.2109896570583666

$ws.Range("B74").Value = "week"
$ws.Range("E74").Value = 0.180252776571334

$ws.Range("B75").Value = "year"
$ws.Range("C75").Value = "arrivaldelay"
$ws.Range("D75").Value = "avg"
$ws.Range("E75").Value = 0.1659470037561426

$ws.Range("B76").Value = "month"
$ws.Range("D76").Value = "avg"
$ws.Range("E76").Value = 0.1658655297543124

$ws.Range("E77").Value = 0.1641404414230867

$ws.Range("B78").Value = "day"
$ws.Range("E78").Value = 0.1626765019637477

$ws.Range("E79").Value = 0.1323626201567746

$ws.Range("B80").Value = "month"
$ws.Range("C80").Value = "arrivaldelay"
$ws.Range("E80").Value = 0.1236558623787677

$ws.Range("B81").Value = "week"
$ws.Range("E81").Value = 0.07120545903503188
